$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'281.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.03%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.73%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.935"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'5.03%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06413"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.39%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.86%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.349"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.65%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8859"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.06%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'1.040"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'14.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1497"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.22%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.05152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.00%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.11%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03100"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.56%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09049"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001571"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.89%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006309"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'3.94%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006044"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.75%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'1.62%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'5.72%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.76%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1290"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.68%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.951"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-3.26%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04357"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.28%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001174"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.40%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.003690"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-9.16%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001197"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.28%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001691"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.54%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04101"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.55%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006646"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1180"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.76%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002355"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'11.80%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01309"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'13.97%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005250"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.78%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E47").Value = "'815.85%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.02246"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-8.25%"
$ws.Range("E48").Style = "Normal"
